# Fruta / hortaliza, semanal
# Insert a new weekly record above the current row 3 (shifting the
# existing rows 3-35 down to 4-36) and populate it with the new
# Guayaba / Vega Modelo de Temuco price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data rows down by one to make room for the new record.
$ws.Rows("3:3").Insert()

# Fill in the newly inserted row 3 with the new weekly observation.
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 45092
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100108
$ws.Range("H3").Value = "Tropicales y subtropicales"
$ws.Range("I3").Value = 100108001
$ws.Range("J3").Value = "Guayaba"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 2600
$ws.Range("O3").Value = 2600
$ws.Range("P3").Value = 2600
$ws.Range("Q3").Value = "$/kilo"
$ws.Range("R3").Value = "Región de Arica y Parinacota"
$ws.Range("S3").Value = 2600
$ws.Range("T3").Value = 1
